$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.216628313064575
$ws.Range("B1").Value = 2.456719398498535
$ws.Range("C1").Value = 7.294825077056885
$ws.Range("D1").Value = 2.254568338394165
$ws.Range("E1").Value = 1.162895202636719
